$wb = $excel.ActiveWorkbook

# Update the "想去人数" (F column) counts on both the "展览" and "全部类型" sheets
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 374
    $ws.Range("F3").Value = 1269
    $ws.Range("F4").Value = 1563
    $ws.Range("F6").Value = 6176
}
